# The commit adds one new weekly price record (row) to the "Tuna" price
# sheet. The new record is inserted as row 8, which pushes the previous
# rows 8-34 down to rows 9-35 (dimension grows from A1:T34 to A1:T35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8, shifting existing rows 8-34 down
# to 9-35.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(8, 3).Value = "Maule"
$ws.Cells.Item(8, 4).Value = 44630
$ws.Cells.Item(8, 5).Value = 7
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100107
$ws.Cells.Item(8, 8).Value = "Otros"
$ws.Cells.Item(8, 9).Value = 100107011
$ws.Cells.Item(8, 10).Value = "Tuna"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Especial"
$ws.Cells.Item(8, 13).Value = 150
$ws.Cells.Item(8, 14).Value = 20000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 20000
$ws.Cells.Item(8, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 20
